# Update Jogos_da_Semana_FlashScore_2025-04-29.xlsx with refreshed FlashScore odds
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("N6").Value = 1.84
$ws.Range("O6").Value = 1.89

# Row 7
$ws.Range("G7").Value = 2.5
$ws.Range("I7").Value = 2.82
$ws.Range("J7").Value = 1.1
$ws.Range("L7").Value = 1.47
$ws.Range("M7").Value = 2.32
$ws.Range("N7").Value = 2.35
$ws.Range("O7").Value = 1.47
$ws.Range("P7").Value = 1.55
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1.65
$ws.Range("T7").Value = 6.4
$ws.Range("U7").Value = 11
$ws.Range("V7").Value = 10.25
$ws.Range("W7").Value = 27
$ws.Range("X7").Value = 25
$ws.Range("Y7").Value = 45
$ws.Range("Z7").Value = 6.7
$ws.Range("AB7").Value = 18
$ws.Range("AC7").Value = 110
$ws.Range("AE7").Value = 6.8
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 35
$ws.Range("AI7").Value = 30
$ws.Range("AJ7").Value = 50

# Row 8
$ws.Range("G8").Value = 2.07
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 3.65
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.57
$ws.Range("T8").Value = 5.5
$ws.Range("U8").Value = 8.5
$ws.Range("V8").Value = 9.25
$ws.Range("W8").Value = 19
$ws.Range("AC8").Value = 150
$ws.Range("AE8").Value = 7.7
$ws.Range("AF8").Value = 17.5
$ws.Range("AG8").Value = 14

# Row 9
$ws.Range("G9").Value = 2.85
$ws.Range("H9").Value = 2.4
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 1.18
$ws.Range("K9").Value = 4.25
$ws.Range("N9").Value = 3
$ws.Range("O9").Value = 1.34
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 2.18
$ws.Range("S9").Value = 1.62
$ws.Range("V9").Value = 11
$ws.Range("Z9").Value = 4.25
$ws.Range("AB9").Value = 17.5
$ws.Range("AE9").Value = 6.4
$ws.Range("AF9").Value = 15

# Row 10
$ws.Range("G10").Value = 2.2
$ws.Range("J10").Value = 1.05
$ws.Range("L10").Value = 1.37
$ws.Range("O10").Value = 1.6
$ws.Range("S10").Value = 1.75

# Row 13
$ws.Range("G13").Value = 1.55
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 5.2
$ws.Range("L13").Value = 1.29
$ws.Range("M13").Value = 3.3
$ws.Range("N13").Value = 1.8
$ws.Range("O13").Value = 1.8
$ws.Range("R13").Value = 1.97
$ws.Range("S13").Value = 1.75
$ws.Range("T13").Value = 5.6
$ws.Range("U13").Value = 6
$ws.Range("V13").Value = 6.9
$ws.Range("X13").Value = 10.5
$ws.Range("Y13").Value = 22
$ws.Range("Z13").Value = 9.75
$ws.Range("AA13").Value = 6.4
$ws.Range("AB13").Value = 14
$ws.Range("AC13").Value = 60
$ws.Range("AD13").Value = 450
$ws.Range("AE13").Value = 11.5
$ws.Range("AF13").Value = 25
$ws.Range("AG13").Value = 14
$ws.Range("AH13").Value = 75
$ws.Range("AJ13").Value = 45

# Row 14
$ws.Range("J14").Value = 1.07
$ws.Range("K14").Value = 9
$ws.Range("L14").Value = 1.36
$ws.Range("O14").Value = 1.67

# Row 15
$ws.Range("J15").Value = 1.1
$ws.Range("L15").Value = 1.5
$ws.Range("O15").Value = 1.5

# Row 16
$ws.Range("G16").Value = 2.55
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 1.11
$ws.Range("M16").Value = 2.5
$ws.Range("N16").Value = 2.6
$ws.Range("O16").Value = 1.48
$ws.Range("Z16").Value = 6.5
$ws.Range("AC16").Value = 67

# Row 17
$ws.Range("J17").Value = 1.07
$ws.Range("L17").Value = 1.33
$ws.Range("O17").Value = 1.7

# Row 18
$ws.Range("G18").Value = 2.75
$ws.Range("I18").Value = 2.63
$ws.Range("J18").Value = 1.08
$ws.Range("K18").Value = 7.5
$ws.Range("L18").Value = 1.4
$ws.Range("M18").Value = 2.75
$ws.Range("N18").Value = 2.35
$ws.Range("O18").Value = 1.57
$ws.Range("P18").Value = 1.5
$ws.Range("Q18").Value = 2.5
$ws.Range("R18").Value = 2
$ws.Range("S18").Value = 1.75
$ws.Range("T18").Value = 7.5
$ws.Range("Z18").Value = 7.5
$ws.Range("AC18").Value = 51
$ws.Range("AD18").Value = 401

# Row 19
$ws.Range("O19").Value = 1.7

# Row 20
$ws.Range("N20").Value = 2.08
$ws.Range("O20").Value = 1.73

# Row 21
$ws.Range("O21").Value = 1.5

# Row 22
$ws.Range("L22").Value = 1.53
$ws.Range("M22").Value = 2.38
$ws.Range("N22").Value = 2.7

# Row 23
$ws.Range("L23").Value = 1.2
$ws.Range("M23").Value = 4.33
$ws.Range("N23").Value = 1.65
$ws.Range("O23").Value = 2.2
$ws.Range("P23").Value = 1.3
$ws.Range("Q23").Value = 3.4
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.25
$ws.Range("Y23").Value = 26
$ws.Range("AE23").Value = 10

# Row 27
$ws.Range("G27").Value = 2.8
$ws.Range("I27").Value = 2.35
$ws.Range("L27").Value = 1.33
$ws.Range("M27").Value = 3.25
$ws.Range("N27").Value = 2.05
$ws.Range("O27").Value = 1.75
$ws.Range("R27").Value = 1.91
$ws.Range("S27").Value = 1.8
$ws.Range("AD27").Value = 351
$ws.Range("AE27").Value = 7

# Row 28
$ws.Range("J28").Value = 1.08
$ws.Range("K28").Value = 8
$ws.Range("L28").Value = 1.44
$ws.Range("M28").Value = 2.63
$ws.Range("N28").Value = 2.4
$ws.Range("O28").Value = 1.53

# Row 29
$ws.Range("N29").Value = 2.05
$ws.Range("O29").Value = 1.8

# Row 34
$ws.Range("G34").Value = 3.2
$ws.Range("I34").Value = 2.27
$ws.Range("L34").Value = 1.5
$ws.Range("M34").Value = 2.47
$ws.Range("N34").Value = 2.47
$ws.Range("P34").Value = 1.55
$ws.Range("Q34").Value = 2.37
$ws.Range("R34").Value = 2.12
$ws.Range("T34").Value = 7.4
$ws.Range("U34").Value = 15.5
$ws.Range("V34").Value = 13
$ws.Range("W34").Value = 45
$ws.Range("X34").Value = 40
$ws.Range("Y34").Value = 60
$ws.Range("AE34").Value = 5.9
$ws.Range("AF34").Value = 10
$ws.Range("AG34").Value = 10.5
$ws.Range("AH34").Value = 23
$ws.Range("AI34").Value = 25
$ws.Range("AJ34").Value = 50

# Row 36
$ws.Range("G36").Value = 1.93
$ws.Range("H36").Value = 3.3
$ws.Range("L36").Value = 1.32
$ws.Range("M36").Value = 2.85
$ws.Range("N36").Value = 1.93
$ws.Range("O36").Value = 1.7
$ws.Range("R36").Value = 1.78
$ws.Range("S36").Value = 1.82
$ws.Range("T36").Value = 6.9
$ws.Range("U36").Value = 9
$ws.Range("V36").Value = 8.5
$ws.Range("X36").Value = 16
$ws.Range("Y36").Value = 29
$ws.Range("Z36").Value = 9
$ws.Range("AA36").Value = 6.4
$ws.Range("AB36").Value = 15.5
$ws.Range("AC36").Value = 75
$ws.Range("AD36").Value = 700
$ws.Range("AE36").Value = 9.75
$ws.Range("AG36").Value = 12.5
$ws.Range("AJ36").Value = 45

# Row 38
$ws.Range("J38").Value = ""
$ws.Range("K38").Value = ""
$ws.Range("L38").Value = 1.05

# Row 39
$ws.Range("J39").Value = 1.02
$ws.Range("L39").Value = 1.13

# Row 40
$ws.Range("G40").Value = 1.9
$ws.Range("H40").Value = 3.4
$ws.Range("I40").Value = 3.8
$ws.Range("J40").Value = 1.03
$ws.Range("L40").Value = 1.25
$ws.Range("P40").Value = 1.36
$ws.Range("Q40").Value = 3
$ws.Range("U40").Value = 9
$ws.Range("X40").Value = 15
$ws.Range("AF40").Value = 21
$ws.Range("AI40").Value = 34
